$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.644726333333334
$ws.Range("H2").Value = 4.934179
$ws.Range("I2").Value = 0.03084360558270512
$ws.Range("J2").Value = 0.03084360558270512
$ws.Range("Q2").Value = 0.6555774998665557
$ws.Range("R2").Value = 5.900197498799001
$ws.Range("S2").Value = 0.03084360558270512
$ws.Range("T2").Value = 0.03084360558270512

# Row 3 updates
$ws.Range("I3").Value = 0.828024694817689
$ws.Range("J3").Value = 0.828024694817689
$ws.Range("S3").Value = 0.828024694817689
$ws.Range("T3").Value = 0.828024694817689

# Row 4 updates
$ws.Range("I4").Value = 0.1411316995996059
$ws.Range("J4").Value = 0.1411316995996059
$ws.Range("S4").Value = 0.1411316995996059
$ws.Range("T4").Value = 0.1411316995996059
